# Completed the question with annotations
# Fills in the Constants sheet with the remaining configuration rows
# (input/output file+sheet names, browser path, timeout, exception
# messages) and adds descriptions for the two input rows that were
# previously missing them.

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# Constants sheet - fix up existing rows 18/19 and populate new rows
# 20-25 with the Unicorn-names exercise configuration.
# ---------------------------------------------------------------------

# Row 18: InputDataFile - path made relative to the workbook + description added
$constants.Cells.Item(18, 2).Value = "Data\Input\input_Unicorn_names.xlsx"
$constants.Cells.Item(18, 3).Value = "Name of the input file"

# Row 19: InputDataSheet - description added
$constants.Cells.Item(19, 3).Value = "Name o f the input sheet"

# Row 20: BrowserPath
$constants.Cells.Item(20, 1).Value = "BrowserPath"
$constants.Cells.Item(20, 2).Value = "http://www.rpasamples.com/unicornname"
$constants.Cells.Item(20, 2).Interior.Pattern = 1
$constants.Cells.Item(20, 2).Interior.Pattern = -4142
$constants.Cells.Item(20, 2).Value = "http://www.rpasamples.com/unicornname"
$constants.Cells.Item(20, 3).Value = "Unicorn browser path"

# Row 21: TimeOut
$constants.Cells.Item(21, 1).Value = "TimeOut"
$constants.Cells.Item(21, 2).Value = 5
$constants.Cells.Item(21, 3).Value = "Delay period for check app state"

# Row 22: OutputDataFile
$constants.Cells.Item(22, 1).Value = "OutputDataFile"
$constants.Cells.Item(22, 2).Value = "Data\Output\Question1.xlsx"
$constants.Cells.Item(22, 3).Value = "Name of the output file"

# Row 23: OutputDataSheet
$constants.Cells.Item(23, 1).Value = "OutputDataSheet"
$constants.Cells.Item(23, 2).Value = "Sheet2"
$constants.Cells.Item(23, 3).Value = "Name of the output sheet"

# Row 24: SystemException
$constants.Cells.Item(24, 1).Value = "SystemException"
$constants.Cells.Item(24, 2).Value = "Page not loading"
$constants.Cells.Item(24, 3).Value = "If the system exception occurs , provide this message"

# Row 25: BusinessException
$constants.Cells.Item(25, 1).Value = "BusinessException"
$constants.Cells.Item(25, 2).Value = "Name and month is incorrect"
$constants.Cells.Item(25, 3).Value = "If the business exception occurs , provide this message"

# Row heights that Excel recalculated for wrapped description cells
$constants.Rows.Item(2).RowHeight = 30
$constants.Rows.Item(3).RowHeight = 45
$constants.Rows.Item(17).RowHeight = 45

# ---------------------------------------------------------------------
# Settings sheet - matching row-height refresh
# ---------------------------------------------------------------------
$settings.Rows.Item(3).RowHeight = 45
$settings.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved view state
# ---------------------------------------------------------------------
$settings.Range("A2").Select()
$constants.Activate()
$constants.Range("B25").Select()
